$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- AIMS box (Content Placeholder 7): fix "Lidar" -> "LiDAR" ---
$aims = $s.Shapes.Item("Content Placeholder 7").TextFrame.TextRange
$aims.Paragraphs(5).Runs(1).Text = "TurtleBot utilizes the LiDAR to enhance position detection"

# --- RESULTS box (Content Placeholder 8): punctuation tweaks, a trimmed
#     sentence, and one brand-new bullet at the end. ---
$results = $s.Shapes.Item("Content Placeholder 8").TextFrame.TextRange
$results.Paragraphs(2).Runs(1).Text = "TurtleBot can be controlled via the ultrasonic sensors."
$results.Paragraphs(3).Runs(1).Text = "Ultrasonic sensors can be used to detect gestures."
$results.Paragraphs(4).Runs(1).Text = "The m5core2 can communicate to the TurtleBot via MQTT."

$lastResults = $results.Paragraphs(5)
$lastResults.Runs(1).Text = "The m5core2 can display the TurtleBot’s position on its display."
[void]$lastResults.InsertAfter("`rThe LiDAR can be utilized to track the TurtleBot’s position.")

# --- CONCLUSIONS box (Content Placeholder 9): each bullet's text shifts
#     down one slot, with a new final "LiDAR" conclusion. ---
$conclusions = $s.Shapes.Item("Content Placeholder 9").TextFrame.TextRange
$conclusions.Paragraphs(3).Runs(1).Text = "Achieved the TurtleBot’s position can be tracked and displayed on the m5core2"
$conclusions.Paragraphs(4).Runs(1).Text = "Achieved that the TurtleBot can communicate with the m5core2 and display its position"
$conclusions.Paragraphs(5).Runs(1).Text = "Achieved that the TurtleBot utilizes the LiDAR to enhance position detection"
